# Updates the histogram sample data: A1/B1 become numeric (-1, 0) instead of
# the "p"/"y" header labels, and the bin counts in column B (B2:B41) are
# refreshed with new values. The bar chart on the sheet reads its cached
# values from Sheet1!B1:B41, so it will pick up the new numbers once Excel
# recalculates.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 1: replace the "p"/"y" text labels with numeric values
$ws.Range("A1").Value = -1
$ws.Range("B1").Value = 0

# New bin counts for B2:B41 (row 42 keeps its SUM(B1:B41) formula)
$newValues = @(
    62,
    147,
    105,
    38,
    18,
    8,
    12,
    16,
    28,
    18,
    8,
    15,
    5,
    17,
    10,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    7,
    7,
    5,
    20,
    3,
    25,
    0,
    9,
    19,
    38,
    37,
    86,
    87,
    103
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $newValues[$i]
}

$wb.Application.Calculate()
